$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.623.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.958.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.58%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.594'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.07'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.365'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0732'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0972'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.251.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.736'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -10.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.981.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.597.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0788'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '219.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -11.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.121'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.115'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.83%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.403.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0200'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0866'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.977'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '14.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.143.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.53%  '
